$wb = $excel.ActiveWorkbook

# --- Sheet: re_profiles ---
$ws4 = $wb.Worksheets.Item("re_profiles")
$ws4.Range("M2").Value = "~TFM_INS-AT"
$ws4.Range("M3").Value = "timeslice"
$ws4.Range("N3").Value = "ncap_afs"
$ws4.Range("O3").Value = "pset_ci"
$ws4.Range("M4").Value = "AllS"
$ws4.Range("N4").Value = 1.2
$ws4.Range("O4").Value = "hydro"

# --- Sheet: load_shapes ---
$ws5 = $wb.Worksheets.Item("load_shapes")
$ws5.Range("I4").Value = "AllSaAllH"
$ws5.Range("J4").Value = 1
$ws5.Range("N4").Value = "AllSaAllH"
$ws5.Range("H5").Value = "elc_industry"
$ws5.Range("I5").Value = "AllSaAllH"
$ws5.Range("J5").Value = 1
$ws5.Range("M5:O5").ClearContents()
$ws5.Rows("6:11").ClearContents()
